# dionisos_Weather_fc_followup.xlsx follow-up update
#
# New monthly "Actual" readings came in for Jul/Aug/Sep/Oct 2021
# (rows 14-17 on Sheet1) and the previously logged value for that
# same month was corrected. Column F ("Diff (fc-act)") is the shared
# formula =D-E, and the two charts plot columns E/F straight out of
# the worksheet, so updating these four cells is enough to ripple the
# forecast-vs-actual delta and the chart series through automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("E14").Value = 22.7   # 2021-07: Actual
$ws.Range("E15").Value = 26.5   # 2021-08: Actual (was blank)
$ws.Range("E16").Value = 27.1   # 2021-09: Actual (was blank)
$ws.Range("E17").Value = 20.4   # 2021-10: Actual (was blank)

$excel.Calculate()
$wb.Save()
